$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 0
$ws1.Range("F11").Value = 320
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 6204
$ws1.Range("F16").Value = 1094
$ws1.Range("F17").Value = 401
$ws1.Range("F18").Value = 31
$ws1.Range("F20").Value = 0
$ws1.Range("F21").Value = 206
$ws1.Range("F24").Value = 10077
$ws1.Range("F27").Value = 1895
$ws1.Range("F33").Value = 165
$ws1.Range("F35").Value = 0
$ws1.Range("F37").Value = 308
$ws1.Range("F38").Value = 61
$ws1.Range("F40").Value = 1197
$ws1.Range("F42").Value = 109
$ws1.Range("F44").Value = 1100
$ws1.Range("F45").Value = 1075
$ws1.Range("F46").Value = 983
$ws1.Range("F48").Value = 63

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 15
$ws2.Range("F7").Value = 31
$ws2.Range("F9").Value = 42
$ws2.Range("F10").Value = 0
$ws2.Range("F11").Value = 1
$ws2.Range("F12").Value = 9
$ws2.Range("F13").Value = 2
$ws2.Range("F15").Value = 0
$ws2.Range("F17").Value = 7
$ws2.Range("F18").Value = 904
$ws2.Range("F19").Value = 0

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F9").Value = 146
$ws4.Range("F10").Value = 241
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 184
$ws4.Range("F16").Value = 5447
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 401
$ws4.Range("F22").Value = 272
$ws4.Range("F23").Value = 147
$ws4.Range("F26").Value = 192
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 0
$ws4.Range("F29").Value = 1916
$ws4.Range("F30").Value = 1895
$ws4.Range("F31").Value = 47
$ws4.Range("F32").Value = 2092
$ws4.Range("F34").Value = 87
$ws4.Range("F35").Value = 0
$ws4.Range("F36").Value = 15
$ws4.Range("F39").Value = 308
$ws4.Range("F40").Value = 5192
$ws4.Range("F41").Value = 0
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 109
$ws4.Range("F44").Value = 0
$ws4.Range("F46").Value = 1075
$ws4.Range("F47").Value = 983
$ws4.Range("F48").Value = 1365
$ws4.Range("F49").Value = 63
$ws4.Range("F50").Value = 1094
